$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.915.64'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.76%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.887.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.20%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.28'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.74%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4617'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.51%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.57'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.05%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9903'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.68'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.35%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.869.40'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.909'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.062'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.70%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.97'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.34%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06554'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.82%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.45'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.37%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.914.54'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.396'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.30%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.22'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.208'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.105.53'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.67'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.63'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.114'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.398'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.83'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.17%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9755'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09350'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.412'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.601'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.27%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.273'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06050'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.48%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02227'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.252'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.60%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.183'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.48%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5770'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.11'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.62%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1818'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.71%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.264'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.283'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +9.64%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.00'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.37%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5456'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.902'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07009'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -7.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.53'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.21%  '
